# création vue initialisation projet
# Shift all dates in column A (rows 3-63) from 2015xxxx to 2017xxxx (+20000),
# and update the corresponding "nombre d'heures" values in column E for the
# rows where they changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every date value in column A (rows 3 through 63) by +20000
for ($r = 3; $r -le 63; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value2 = $current + 20000
    }
}

# New values for column E (only rows whose value actually changed)
$eUpdates = @{
    3  = 12
    4  = 20
    5  = 20
    6  = 20
    7  = 15
    9  = 20
    10 = 7
    11 = 11
    12 = 15
    13 = 18
    14 = 8
    15 = 13
    16 = 12
    18 = 13
    19 = 8
    20 = 13
    21 = 5
    22 = 11
    23 = 16
    24 = 18
    25 = 5
    26 = 20
    27 = 7
    28 = 5
    29 = 9
    30 = 9
    31 = 9
    32 = 13
    33 = 14
    34 = 9
    35 = 14
    36 = 17
    37 = 17
    38 = 5
    39 = 19
    40 = 20
    41 = 10
    42 = 6
    43 = 7
    44 = 6
    45 = 19
    46 = 13
    47 = 9
    48 = 11
    49 = 8
    50 = 9
    52 = 6
    53 = 20
    54 = 12
    55 = 18
    56 = 16
    57 = 7
    58 = 11
    59 = 7
    60 = 5
    61 = 16
    62 = 11
    63 = 11
}

foreach ($row in $eUpdates.Keys) {
    $ws.Cells.Item([int]$row, 5).Value2 = $eUpdates[$row]
}
